# "improved and cut some" - the "classical" genre was dropped from the dataset
# and replaced with "soul" everywhere it appeared. In the two rows where that
# substitution would have produced a duplicate genre within the same row
# (row 18 and row 42, which already contained "soul"), the colliding cell was
# instead swapped for a different genre so every row keeps three distinct
# genres.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value
$updates = @{
    "E3"  = "soul"
    "F5"  = "soul"
    "D6"  = "soul"
    "D11" = "soul"
    "F14" = "soul"
    "E18" = "soul"
    "F18" = "country"
    "F23" = "soul"
    "F34" = "soul"
    "E35" = "soul"
    "F36" = "soul"
    "F41" = "soul"
    "E42" = "soul"
    "F42" = "folk"
    "E49" = "soul"
    "E51" = "soul"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# The Name column (C) ends up auto-fit to its contents after the edit.
$ws.Columns.Item(3).AutoFit()

# Leave the same cell selected as in the saved workbook.
$ws.Range("D5").Select()
